$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New "Multiplexer Select" rows under the "Other" table (B21:B23) ---
$ws.Range("B21").Value = "Multiplexer Select 1"
$ws.Range("B22").Value = "Multiplexer Select 2"
$ws.Range("B23").Value = "Multiplexer Select 3"

# --- 4051 multiplexer pinout diagram (columns G:J, rows 13-21) ---
$ws.Range("H13").Value = "4051 Pinout"

$ws.Range("G14").Value = "CH 4 I/O"
$ws.Range("H14").Value = 1
$ws.Range("I14").Value = 16
$ws.Range("J14").Value = "VDD"
$ws.Range("J14").Characters(2,2).Font.Size = 8
$ws.Range("J14").Characters(2,2).Font.Bold = $false
$ws.Range("J14").Characters(2,2).Font.ThemeColor = 1

$ws.Range("G15").Value = "CH 6 I/O"
$ws.Range("H15").Value = 2
$ws.Range("I15").Value = 15
$ws.Range("J15").Value = "CH 2 I/O"

$ws.Range("G16").Value = "COM O/I"
$ws.Range("H16").Value = 3
$ws.Range("I16").Value = 14
$ws.Range("J16").Value = "CH 1 I/O"

$ws.Range("G17").Value = "CH 7 I/O"
$ws.Range("H17").Value = 4
$ws.Range("I17").Value = 13
$ws.Range("J17").Value = "CH 0 I/O"

$ws.Range("G18").Value = "CH 5 I/O"
$ws.Range("H18").Value = 5
$ws.Range("I18").Value = 12
$ws.Range("J18").Value = "CH 3 I/O"

$ws.Range("G19").Value = "INH"
$ws.Range("H19").Value = 6
$ws.Range("I19").Value = 11
$ws.Range("J19").Value = "A"

$ws.Range("G20").Value = "VEE"
$ws.Range("G20").Characters(2,2).Font.Size = 8
$ws.Range("G20").Characters(2,2).Font.Bold = $false
$ws.Range("G20").Characters(2,2).Font.ThemeColor = 1
$ws.Range("H20").Value = 7
$ws.Range("I20").Value = 10
$ws.Range("J20").Value = "B"

$ws.Range("G21").Value = "VSS"
$ws.Range("G21").Characters(2,2).Font.Size = 8
$ws.Range("G21").Characters(2,2).Font.Bold = $false
$ws.Range("G21").Characters(2,2).Font.ThemeColor = 1
$ws.Range("H21").Value = 8
$ws.Range("I21").Value = 9
$ws.Range("J21").Value = "C"

# --- Formatting for the pinout number columns (H and I, rows 14-21) ---
$hRange = $ws.Range("H14:H21")
$hRange.Interior.Pattern = -4124
$hRange.Interior.ThemeColor = 1
$hRange.Interior.TintAndShade = 0.249977111117893
$hRange.NumberFormat = "0"
$hRange.HorizontalAlignment = -4131

$iRange = $ws.Range("I14:I21")
$iRange.Interior.Pattern = -4124
$iRange.Interior.ThemeColor = 1
$iRange.Interior.TintAndShade = 0.249977111117893

# --- View adjustments ---
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("J25").Select()
